$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting existing rows 93:150 down to 94:151
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new data record
$ws.Cells.Item(93, 1).Value = 5
$ws.Cells.Item(93, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(93, 3).Value = "Maule"
$ws.Cells.Item(93, 4).Value = 45233
$ws.Cells.Item(93, 5).Value = 7
$ws.Cells.Item(93, 6).Value = 100112022
$ws.Cells.Item(93, 7).Value = "Arveja Verde"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 500
$ws.Cells.Item(93, 11).Value = 18000
$ws.Cells.Item(93, 12).Value = 20000
$ws.Cells.Item(93, 13).Value = 18800
$ws.Cells.Item(93, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(93, 15).Value = "Región del Maule"
$ws.Cells.Item(93, 16).Value = 752
$ws.Cells.Item(93, 17).Value = 25
$ws.Cells.Item(93, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date-number format style as other rows in column D
$ws.Range("D93").NumberFormat = $ws.Range("D94").NumberFormat
